$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 211 (current data row 211..271 shifts down to 212..272),
# mirroring the new weekly observation added at the top of this block.
$ws.Rows("211:211").Insert()

# Populate the freshly inserted row 211 with the new weekly record.
$ws.Range("A211").Value = 3
$ws.Range("B211").Value = "Femacal de La Calera"
$ws.Range("C211").Value = "Coquimbo"
$ws.Range("D211").Value = 45120
$ws.Range("E211").Value = 5
$ws.Range("F211").Value = 100112026
$ws.Range("G211").Value = "Haba"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 45
$ws.Range("K211").Value = 18000
$ws.Range("L211").Value = 18000
$ws.Range("M211").Value = 18000
$ws.Range("N211").Value = "$/saco 25 kilos"
$ws.Range("O211").Value = "Provincia de Limarí"
$ws.Range("P211").Value = 720
$ws.Range("Q211").Value = 25
$ws.Range("R211").Value = "Hortaliza"
